$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 306: the "close" (F) value was revised from 30.75 to 30.85
$ws.Range("F306").Value = 30.85

# Copy row 306's formatting (bold/date-format/border on col A, borders on B:G)
# down into the three new rows so they match the sheet's existing style.
$ws.Range("A306:G306").Copy()
$ws.Range("A307:G309").PasteSpecial(-4122)

# New row 307 - 2023-05-01
$ws.Range("A307").Value = 45047.33333333334
$ws.Range("B307").Value = "FX_IDC:USDEGP"
$ws.Range("C307").Value = 30.9499
$ws.Range("D307").Value = 30.9499
$ws.Range("E307").Value = 30.73
$ws.Range("F307").Value = 30.85
$ws.Range("G307").Value = 0

# New row 308 - 2023-06-01
$ws.Range("A308").Value = 45078.33333333334
$ws.Range("B308").Value = "FX_IDC:USDEGP"
$ws.Range("C308").Value = 30.85
$ws.Range("D308").Value = 30.9499
$ws.Range("E308").Value = 30.75
$ws.Range("F308").Value = 30.85
$ws.Range("G308").Value = 0

# New row 309 - 2023-07-03
$ws.Range("A309").Value = 45110.33333333334
$ws.Range("B309").Value = "FX_IDC:USDEGP"
$ws.Range("C309").Value = 30.85
$ws.Range("D309").Value = 30.9499
$ws.Range("E309").Value = 30.75
$ws.Range("F309").Value = 30.83
$ws.Range("G309").Value = 0
